# Y5_B2526_General_&_Special_Surgery_2_B1_schedule.xlsx
#
# For every group's (B1-1 .. B1-12) 08:30 sub-block (sessions 16-22, i.e. the
# last 7 rows of each 22-row group block), the "Subject" column (C) changes
# from "general surgery" to the new value "Surgery Seminar/Slide". Each
# 22-row group block starts at rows 17, 39, 61, 83, 105, 127, 149, 171, 193,
# 215, 237, 259; the affected sub-block is the 7 rows starting at each of
# those rows.
#
# Also bump the sheet's zoom level to 156%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSubject = "Surgery Seminar/Slide"

# Starting row of each group's 22-row block.
$blockStarts = @(17, 39, 61, 83, 105, 127, 149, 171, 193, 215, 237, 259)

# Update the cell values for all 7 rows (offsets 0-6) in every block.
foreach ($blockStart in $blockStarts) {
    for ($offset = 0; $offset -le 6; $offset++) {
        $row = $blockStart + $offset
        $ws.Cells.Item($row, 3).Value = $newSubject
    }
}

# Re-apply the "unshaded" formatting (matching the sibling rows in the same
# block) to the rows that previously had the shaded style (offsets 1, 3, 5
# -> style index 2) so they match the rest of the block (style index 6), by
# copying the format from the first (already correctly-styled) row of each
# block.
foreach ($blockStart in $blockStarts) {
    $ws.Cells.Item($blockStart, 3).Copy() | Out-Null
    for ($offset = 1; $offset -le 5; $offset += 2) {
        $row = $blockStart + $offset
        $ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null
    }
}
$excel.CutCopyMode = 0

# Bump the zoom level of the active sheet's view to 156%.
$excel.ActiveWindow.Zoom = 156
